$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row with new column names
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# Title-case the Spanish connector words (de, del, el, la, las, los, y)
# in state/municipality names throughout the data
$ws.Range("B6").Value = 'Pabellón De Arteaga'
$ws.Range("B7").Value = 'Rincón De Romos'
$ws.Range("B8").Value = 'San Francisco De Los Romo'
$ws.Range("B9").Value = 'San José De Gracia'
$ws.Range("B31").Value = 'Comitán De Domínguez'
$ws.Range("B49").Value = 'Salto De Agua'
$ws.Range("B50").Value = 'San Cristóbal De Las Casas'
$ws.Range("B73").Value = 'Guadalupe Y Calvo'
$ws.Range("B74").Value = 'Hidalgo Del Parral'
$ws.Range("B84").Value = 'San Francisco De Borja'
$ws.Range("B87").Value = 'Valle De Zaragoza'
$ws.Range("B101").Value = 'San Juan De Sabinas'
$ws.Range("B110").Value = 'Villa De Álvarez'
$ws.Range("A112").Value = 'Ciudad De México'
$ws.Range("B116").Value = 'Cuajimalpa De Morelos'
$ws.Range("B140").Value = 'Nombre De Dios'
$ws.Range("A156").Value = 'Estado De México'
$ws.Range("B156").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B159").Value = 'Almoloya De Juárez'
$ws.Range("B160").Value = 'Almoloya Del Río'
$ws.Range("B166").Value = 'Atizapán De Zaragoza'
$ws.Range("B176").Value = 'Ecatepec De Morelos'
$ws.Range("B181").Value = 'Ixtapan De La Sal'
$ws.Range("B182").Value = 'Ixtapan Del Oro'
$ws.Range("B192").Value = 'Naucalpan De Juárez'
$ws.Range("B197").Value = 'San Felipe Del Progreso'
$ws.Range("B198").Value = 'San Simón De Guerrero'
$ws.Range("B207").Value = 'Tenango Del Valle'
$ws.Range("B215").Value = 'Tlalnepantla De Baz'
$ws.Range("B221").Value = 'Valle De Bravo'
$ws.Range("B222").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B223").Value = 'Villa De Allende'
$ws.Range("B234").Value = 'San Miguel De Allende'
$ws.Range("B235").Value = 'Apaseo El Alto'
$ws.Range("B236").Value = 'Apaseo El Grande'
$ws.Range("B243").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B246").Value = 'Jaral Del Progreso'
$ws.Range("B253").Value = 'Purísima Del Rincón'
$ws.Range("B257").Value = 'San Diego De La Unión'
$ws.Range("B259").Value = 'San Francisco Del Rincón'
$ws.Range("B261").Value = 'San Luis De La Paz'
$ws.Range("B262").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B264").Value = 'Silao De La Victoria'
$ws.Range("B268").Value = 'Valle De Santiago'
$ws.Range("B274").Value = 'Acapulco De Juárez'
$ws.Range("B276").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B277").Value = 'Alcozauca De Guerrero'
$ws.Range("B280").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B282").Value = 'Atoyac De Álvarez'
$ws.Range("B285").Value = 'Buenavista De Cuéllar'
$ws.Range("B286").Value = 'Chilapa De Álvarez'
$ws.Range("B287").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B288").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B292").Value = 'Coyuca De Benítez'
$ws.Range("B293").Value = 'Coyuca De Catalán'
$ws.Range("B296").Value = 'Cuetzala Del Progreso'
$ws.Range("B297").Value = 'Cutzamala De Pinzón'
$ws.Range("B302").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B303").Value = 'Iguala De La Independencia'
$ws.Range("B304").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B305").Value = 'Zihuatanejo De Azueta'
$ws.Range("B306").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B319").Value = 'Taxco De Alarcón'
$ws.Range("B321").Value = 'Técpan De Galeana'
$ws.Range("B323").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B325").Value = 'Tixtla De Guerrero'
$ws.Range("B328").Value = 'Tlapa De Comonfort'
$ws.Range("B339").Value = 'Cuautepec De Hinojosa'
$ws.Range("B341").Value = 'Huasca De Ocampo'
$ws.Range("B343").Value = 'Huejutla De Reyes'
$ws.Range("B346").Value = 'Jacala De Ledezma'
$ws.Range("B349").Value = 'Molango De Escamilla'
$ws.Range("B351").Value = 'Pachuca De Soto'
$ws.Range("B352").Value = 'Progreso De Obregón'
$ws.Range("B355").Value = 'Santiago De Anaya'
$ws.Range("B359").Value = 'Tepehuacán De Guerrero'
$ws.Range("B360").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B361").Value = 'Tezontepec De Aldama'
$ws.Range("B367").Value = 'Tula De Allende'
$ws.Range("B368").Value = 'Tulancingo De Bravo'
$ws.Range("B371").Value = 'Zacualtipán De Ángeles'
$ws.Range("B372").Value = 'Zapotlán De Juárez'
$ws.Range("B375").Value = 'Ahualulco De Mercado'
$ws.Range("B378").Value = 'Atemajac De Brizuela'
$ws.Range("B379").Value = 'Atotonilco El Alto'
$ws.Range("B380").Value = 'Autlán De Navarro'
$ws.Range("B391").Value = 'Encarnación De Díaz'
$ws.Range("B395").Value = 'Huejuquilla El Alto'
$ws.Range("B396").Value = 'Ixtlahuacán Del Río'
$ws.Range("B399").Value = 'Jilotlán De Los Dolores'
$ws.Range("B403").Value = 'Lagos De Moreno'
$ws.Range("B413").Value = 'San Diego De Alejandría'
$ws.Range("B415").Value = 'San Juan De Los Lagos'
$ws.Range("B417").Value = 'Santa María De Los Ángeles'
$ws.Range("B420").Value = 'Talpa De Allende'
$ws.Range("B421").Value = 'Tamazula De Gordiano'
$ws.Range("B427").Value = 'Tepatitlán De Morelos'
$ws.Range("B429").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B437").Value = 'Unión De San Antonio'
$ws.Range("B438").Value = 'Unión De Tula'
$ws.Range("B441").Value = 'Yahualica De González Gallo'
$ws.Range("B442").Value = 'Zacoalco De Torres'
$ws.Range("B445").Value = 'Zapotlán Del Rey'
$ws.Range("B446").Value = 'Zapotlán El Grande'
$ws.Range("B513").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B534").Value = 'Coatlán Del Río'
$ws.Range("B541").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B545").Value = 'Puente De Ixtla'
$ws.Range("B549").Value = 'Tetela Del Volcán'
$ws.Range("B550").Value = 'Tlaltizapán De Zapata'
$ws.Range("B558").Value = 'Amatlán De Cañas'
$ws.Range("B561").Value = 'Ixtlán Del Río'
$ws.Range("B583").Value = 'San Nicolás De Los Garza'
$ws.Range("B589").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B592").Value = 'Ayoquezco De Aldama'
$ws.Range("B594").Value = 'Capulálpam De Méndez'
$ws.Range("B595").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B596").Value = 'Chiquihuitlán De Benito Juárez'
$ws.Range("B597").Value = 'Coicoyán De Las Flores'
$ws.Range("B598").Value = 'Constancia Del Rosario'
$ws.Range("B600").Value = 'Cuilápam De Guerrero'
$ws.Range("B601").Value = 'Cuyamecalco Villa De Zaragoza'
$ws.Range("B602").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B603").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B604").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B606").Value = 'Ixtlán De Juárez'
$ws.Range("B607").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B614").Value = 'Mariscala De Juárez'
$ws.Range("B616").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B619").Value = 'Oaxaca De Juárez'
$ws.Range("B620").Value = 'Ocotlán De Morelos'
$ws.Range("B621").Value = 'Putla Villa De Guerrero'
$ws.Range("B641").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B655").Value = 'San Juan Del Estado'
$ws.Range("B656").Value = 'San Juan Del Río'
$ws.Range("B683").Value = 'San Miguel Del Puerto'
$ws.Range("B692").Value = 'San Pedro El Alto'
$ws.Range("B702").Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range("B718").Value = 'Santa Inés Del Monte'
$ws.Range("B729").Value = 'Santa María Jalapa Del Marqués'
$ws.Range("B749").Value = 'Santiago Del Río'
$ws.Range("B771").Value = 'Santo Domingo De Morelos'
$ws.Range("B780").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B781").Value = 'Tanetze De Zaragoza'
$ws.Range("B782").Value = 'Tataltepec De Valdés'
$ws.Range("B783").Value = 'Teotitlán De Flores Magón'
$ws.Range("B784").Value = 'Tepelmeme Villa De Morelos'
$ws.Range("B785").Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range("B786").Value = 'Tlacolula De Matamoros'
$ws.Range("B787").Value = 'Totontepec Villa De Morelos'
$ws.Range("B790").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B791").Value = 'Villa De Tututepec'
$ws.Range("B792").Value = 'Villa De Zaachila'
$ws.Range("B794").Value = 'Villa Sola De Vega'
$ws.Range("B795").Value = 'Villa Talea De Castro'
$ws.Range("B796").Value = 'Zapotitlán Del Río'
$ws.Range("B798").Value = 'Zimatlán De Álvarez'
$ws.Range("B808").Value = 'Ayotoxco De Guerrero'
$ws.Range("B810").Value = 'Chalchicomula De Sesma'
$ws.Range("B823").Value = 'Huehuetlán El Chico'
$ws.Range("B824").Value = 'Huehuetlán El Grande'
$ws.Range("B828").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B830").Value = 'Izúcar De Matamoros'
$ws.Range("B837").Value = 'Los Reyes De Juárez'
$ws.Range("B842").Value = 'Palmar De Bravo'
$ws.Range("B856").Value = 'San Salvador El Seco'
$ws.Range("B857").Value = 'San Salvador El Verde'
$ws.Range("B861").Value = 'Tecali De Herrera'
$ws.Range("B867").Value = 'Tepanco De López'
$ws.Range("B868").Value = 'Tepango De Rodríguez'
$ws.Range("B873").Value = 'Tepexi De Rodríguez'
$ws.Range("B878").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B899").Value = 'Amealco De Bonfil'
$ws.Range("B901").Value = 'Cadereyta De Montes'
$ws.Range("B903").Value = 'Jalpan De Serra'
$ws.Range("B904").Value = 'Landa De Matamoros'
$ws.Range("B906").Value = 'Pinal De Amoles'
$ws.Range("B908").Value = 'San Juan Del Río'
$ws.Range("B917").Value = 'Axtla De Terrazas'
$ws.Range("B921").Value = 'Ciudad Del Maíz'
$ws.Range("B931").Value = 'Santa María Del Río'
$ws.Range("B938").Value = 'Tanquián De Escobedo'
$ws.Range("B941").Value = 'Villa De Arista'
$ws.Range("B942").Value = 'Villa De Ramos'
$ws.Range("B943").Value = 'Villa De Reyes'
$ws.Range("B997").Value = 'Soto La Marina'
$ws.Range("B1004").Value = 'Acuamanala De Miguel Hidalgo'
$ws.Range("B1005").Value = 'Amaxac De Guerrero'
$ws.Range("B1012").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B1014").Value = 'Papalotla De Xicohténcatl'
$ws.Range("B1018").Value = 'Tepetitla De Lardizábal'
$ws.Range("B1035").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B1039").Value = 'Amatlán De Los Reyes'
$ws.Range("B1046").Value = 'Boca Del Río'
$ws.Range("B1048").Value = 'Camarón De Tejeda'
$ws.Range("B1053").Value = 'Cazones De Herrera'
$ws.Range("B1063").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1074").Value = 'Hueyapan De Ocampo'
$ws.Range("B1075").Value = 'Ignacio De La Llave'
$ws.Range("B1077").Value = 'Ixhuatlán De Madero'
$ws.Range("B1078").Value = 'Ixhuatlán Del Café'
$ws.Range("B1079").Value = 'Ixhuatlán Del Sureste'
$ws.Range("B1089").Value = 'Juchique De Ferrer'
$ws.Range("B1092").Value = 'Landero Y Coss'
$ws.Range("B1095").Value = 'Lerdo De Tejada'
$ws.Range("B1099").Value = 'Martínez De La Torre'
$ws.Range("B1114").Value = 'Paso Del Macho'
$ws.Range("B1117").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1124").Value = 'Sayula De Alemán'
$ws.Range("B1127").Value = 'Soledad De Doblado'
$ws.Range("B1130").Value = 'Tatahuicapan De Juárez'
$ws.Range("B1152").Value = 'Vega De Alatorre'
$ws.Range("B1179").Value = 'Jiménez Del Teul'
$ws.Range("B1184").Value = 'Mezquital Del Oro'
$ws.Range("B1188").Value = 'Moyahua De Estrada'
$ws.Range("B1189").Value = 'Nochistlán De Mejía'
$ws.Range("B1190").Value = 'Noria De Ángeles'
$ws.Range("B1201").Value = 'Teúl De González Ortega'
$ws.Range("B1202").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1203").Value = 'Trinidad García De La Cadena'
$ws.Range("B1205").Value = 'Villa De Cos'

# Tiny floating point re-serialization (last-bit) adjustments picked up
# when the source pipeline recomputed these percentages
$ws.Range("D231").Value = 0.09304161497688056
$ws.Range("D448").Value = 0.09157550468027516

# Remove the trailing footnote/metadata rows (1213:1217)
$ws.Range("A1213:D1217").ClearContents()
